$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 (Magnesium chloride unit price): update loading value and
# replace the dependent formulas (G8, I8) with their computed results.
$ws.Range("E8").Value = 0.38
$ws.Range("G8").Value = 0.349
$ws.Range("I8").Value = 0.411
$ws.Range("Q8").Value = 1

# Row 9 (Zinc sulfate unit price): same treatment.
$ws.Range("E9").Value = 0.795
$ws.Range("G9").Value = 0.657
$ws.Range("I9").Value = 0.931

# Update the active selection to match the new cursor position left by
# the edits (rows 8-9 selected).
$ws.Range("A8:XFD9").Select() | Out-Null
